$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O:P. This shifts the former O:U content to Q:W,
# preserving their values/formatting, and grows the used range accordingly.
$ws.Columns("O:P").Insert()

# Rename the (now) M1/N1 headers.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Add headers for the two newly inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Populate the new O/P columns for each data row with the same values as the
# corresponding M/N columns (Detected Predicates Doc Parent/Related).
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 15).Value = $ws.Cells.Item($row, 13).Value2
    $ws.Cells.Item($row, 16).Value = $ws.Cells.Item($row, 14).Value2
}
